$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): insert Corequisites, Concurrent, Recommended columns; shift Terms Typically Offered to G1 ---
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"
$ws.Range("G1").Value = "Terms Typically Offered"

# --- Data rows 2-95 ---
# Row 2
$ws.Range("C2").Value = "Mechanical Engineering student; first quarter of freshman year."
$ws.Range("D2").Value = "NA"
$ws.Range("E2").Value = "ME 163."
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "F "

# Row 3
$ws.Range("C3").Value = "ME 128; Mechanical Engineering student; second quarter of freshman year."
$ws.Range("D3").Value = "NA"
$ws.Range("E3").Value = "IME 145."
$ws.Range("F3").Value = "NA"
$ws.Range("G3").Value = "W "

# Row 4
$ws.Range("C4").Value = "ME 129; Mechanical Engineering student; third quarter of freshman year."
$ws.Range("D4").Value = "IME 146."
$ws.Range("E4").Value = "NA"
$ws.Range("F4").Value = "NA"
$ws.Range("G4").Value = "SP "

# Row 5
$ws.Range("C5").Value = "NA"
$ws.Range("D5").Value = "NA"
$ws.Range("E5").Value = "NA"
$ws.Range("F5").Value = "NA"
$ws.Range("G5").Value = "F,W,SP,SU"

# Row 6
$ws.Range("C6").Value = "NA"
$ws.Range("D6").Value = "NA"
$ws.Range("E6").Value = "ME 128."
$ws.Range("F6").Value = "NA"
$ws.Range("G6").Value = "F"

# Row 7
$ws.Range("C7").Value = "MATH 241 (or concurrently), PHYS 131 or PHYS 141."
$ws.Range("D7").Value = "NA"
$ws.Range("E7").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("G7").Value = "F, W, SP"

# Row 8
$ws.Range("C8").Value = "MATH 241; ME 211 or ARCE 211."
$ws.Range("D8").Value = "NA"
$ws.Range("E8").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("G8").Value = "F, W, SP"

# Row 9
$ws.Range("C9").Value = "NA"
$ws.Range("D9").Value = "NA"
$ws.Range("E9").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("G9").Value = "F, W, SP"

# Row 10
$ws.Range("C10").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("E10").Value = "NA"
$ws.Range("F10").Value = "NA"
$ws.Range("G10").Value = "F, W, SP"

# Row 11
$ws.Range("C11").Value = "Sophomore standing."
$ws.Range("D11").Value = "NA"
$ws.Range("E11").Value = "NA"
$ws.Range("F11").Value = "NA"
$ws.Range("G11").Value = "F, W, SP"

# Row 12
$ws.Range("C12").Value = "Engineering majors."
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("F12").Value = "CHEM 125, ENGL 149, and PHYS 132."
$ws.Range("G12").Value = "F, W, SP "

# Row 13
$ws.Range("C13").Value = "ME 130 or ME 228."
$ws.Range("D13").Value = "NA"
$ws.Range("E13").Value = "NA"
$ws.Range("F13").Value = "IME 143."
$ws.Range("G13").Value = "F, W, SP "

# Row 14
$ws.Range("C14").Value = "NA"
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("F14").Value = "NA"
$ws.Range("G14").Value = "F, SP"

# Row 15
$ws.Range("C15").Value = "NA"
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = "NA"
$ws.Range("G15").Value = "F, SP"

# Row 16
$ws.Range("C16").Value = "Open to undergraduate students and consent of instructor."
$ws.Range("D16").Value = "NA"
$ws.Range("E16").Value = "NA"
$ws.Range("F16").Value = "NA"
$ws.Range("G16").Value = "TBD"

# Row 17
$ws.Range("C17").Value = "Consent of instructor."
$ws.Range("D17").Value = "NA"
$ws.Range("E17").Value = "NA"
$ws.Range("F17").Value = "NA"
$ws.Range("G17").Value = "TBD"

# Row 18
$ws.Range("C18").Value = "ME 212 and PHYS 132."
$ws.Range("D18").Value = "NA"
$ws.Range("E18").Value = "NA"
$ws.Range("F18").Value = "NA"
$ws.Range("G18").Value = "F, W, SP"

# Row 19
$ws.Range("C19").Value = "ME 302."
$ws.Range("D19").Value = "NA"
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = "NA"
$ws.Range("G19").Value = "F, W, SP"

# Row 20
$ws.Range("C20").Value = "EE 201 and EE 251."
$ws.Range("D20").Value = "NA"
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = "NA"
$ws.Range("G20").Value = "F, W"

# Row 21
$ws.Range("C21").Value = "ME 212, MATH 344."
$ws.Range("D21").Value = "NA"
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = "EE 201."
$ws.Range("G21").Value = "F, W, SP "

# Row 22
$ws.Range("C22").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; and completion of GE Areas B2, B3, and B4."
$ws.Range("D22").Value = "NA"
$ws.Range("E22").Value = "NA"
$ws.Range("F22").Value = "NA"
$ws.Range("G22").Value = "TBD"

# Row 23
$ws.Range("C23").Value = "NA"
$ws.Range("D23").Value = "NA"
$ws.Range("E23").Value = "NA"
$ws.Range("F23").Value = "NA"
$ws.Range("G23").Value = "F, W, CSC 231 or CSC 234; EE 201; EE 251; ME 318; ME 341."

# Row 24
$ws.Range("C24").Value = "Junior standing; completion of GE Area A with grades of C- or better; completion of GE Area B1 with a grade of C- or better in at least one of the courses; and completion of GE Areas B2, B3, and B4."
$ws.Range("D24").Value = "NA"
$ws.Range("E24").Value = "NA"
$ws.Range("F24").Value = "NA"
$ws.Range("G24").Value = "W"

# Row 25
$ws.Range("C25").Value = "ME 212; CSC 231 or CSC 234."
$ws.Range("D25").Value = "MATH 244."
$ws.Range("E25").Value = "NA"
$ws.Range("F25").Value = "NA"
$ws.Range("G25").Value = "F, W, SP "

# Row 26
$ws.Range("C26").Value = "BMED 212 or ME 234; CE 207; CSC 231 or CSC 234; MATE 210; ME 212; and ME 251."
$ws.Range("D26").Value = "IME 141 or ITP 341."
$ws.Range("E26").Value = "NA"
$ws.Range("F26").Value = "NA"
$ws.Range("G26").Value = "F, W, SP "

# Row 27
$ws.Range("C27").Value = "ME 328."
$ws.Range("D27").Value = "NA"
$ws.Range("E27").Value = "NA"
$ws.Range("F27").Value = "NA"
$ws.Range("G27").Value = "F, W, SP"

# Row 28
$ws.Range("C28").Value = "MATH 242 or MATH 244; ME 212."
$ws.Range("D28").Value = "NA"
$ws.Range("E28").Value = "NA"
$ws.Range("F28").Value = "NA"
$ws.Range("G28").Value = "F, W, SP"

# Row 29
$ws.Range("C29").Value = "ME 236, ME 341, ME 302."
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("F29").Value = "NA"
$ws.Range("G29").Value = "F, W, SP"

# Row 30
$ws.Range("C30").Value = "CPE/CSC 101 or CSC 231 or CSC 234; MATE 360 and MATE 380, or ME 236 and ME 302 and ME 341."
$ws.Range("D30").Value = "NA"
$ws.Range("E30").Value = "NA"
$ws.Range("F30").Value = "NA"
$ws.Range("G30").Value = "F, W, SP"

# Row 31
$ws.Range("C31").Value = "NA"
$ws.Range("D31").Value = "ME 302."
$ws.Range("E31").Value = "NA"
$ws.Range("F31").Value = "NA"
$ws.Range("G31").Value = "W"

# Row 32
$ws.Range("C32").Value = "Consent of instructor."
$ws.Range("D32").Value = "NA"
$ws.Range("E32").Value = "NA"
$ws.Range("F32").Value = "NA"
$ws.Range("G32").Value = "F, W, SP"

# Row 33
$ws.Range("C33").Value = "CE 207, MATH 344, ME 328 or consent of instructor."
$ws.Range("D33").Value = "NA"
$ws.Range("E33").Value = "NA"
$ws.Range("F33").Value = "NA"
$ws.Range("G33").Value = "F"

# Row 34
$ws.Range("C34").Value = "ME 328."
$ws.Range("D34").Value = "NA"
$ws.Range("E34").Value = "NA"
$ws.Range("F34").Value = "NA"
$ws.Range("G34").Value = "SP"

# Row 35
$ws.Range("C35").Value = "BMED 410, and CE 207 or CE 208; or CE 406; or ME 328."
$ws.Range("D35").Value = "NA"
$ws.Range("E35").Value = "NA"
$ws.Range("F35").Value = "NA"
$ws.Range("G35").Value = "F, W, SP"

# Row 36
$ws.Range("C36").Value = "EE 321, EE 361, ME 305, and ME 329 (ME329 may be taken concurrently); CPE 316 or CPE/EE 329 or CPE/EE 336."
$ws.Range("D36").Value = "NA"
$ws.Range("E36").Value = "NA"
$ws.Range("F36").Value = "NA"
$ws.Range("G36").Value = "W, SP"

# Row 37
$ws.Range("C37").Value = "KINE 403; or ME 326."
$ws.Range("D37").Value = "NA"
$ws.Range("E37").Value = "NA"
$ws.Range("F37").Value = "NA"
$ws.Range("G37").Value = "TBD"

# Row 38
$ws.Range("C38").Value = "ME 328."
$ws.Range("D38").Value = "NA"
$ws.Range("E38").Value = "NA"
$ws.Range("F38").Value = "ME 318."
$ws.Range("G38").Value = "SP "

# Row 39
$ws.Range("C39").Value = "AERO 331 or ME 328."
$ws.Range("D39").Value = "NA"
$ws.Range("E39").Value = "NA"
$ws.Range("F39").Value = "NA"
$ws.Range("G39").Value = "F, SP"

# Row 40
$ws.Range("C40").Value = "ME 302."
$ws.Range("D40").Value = "NA"
$ws.Range("E40").Value = "NA"
$ws.Range("F40").Value = "NA"
$ws.Range("G40").Value = "SP"

# Row 41
$ws.Range("C41").Value = "ME 318, ME 328."
$ws.Range("D41").Value = "NA"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "F"

# Row 42
$ws.Range("C42").Value = "ME 322."
$ws.Range("D42").Value = "NA"
$ws.Range("E42").Value = "NA"
$ws.Range("F42").Value = "NA"
$ws.Range("G42").Value = "F, W, SP"

# Row 43
$ws.Range("C43").Value = "ME 322, ME 236."
$ws.Range("D43").Value = "NA"
$ws.Range("E43").Value = "NA"
$ws.Range("F43").Value = "NA"
$ws.Range("G43").Value = "F, W, SP"

# Row 44
$ws.Range("C44").Value = "ME 303; ME 347; and ME 350."
$ws.Range("D44").Value = "NA"
$ws.Range("E44").Value = "NA"
$ws.Range("F44").Value = "NA"
$ws.Range("G44").Value = "F, W, SP"

# Row 45
$ws.Range("C45").Value = "ME 318."
$ws.Range("D45").Value = "NA"
$ws.Range("E45").Value = "NA"
$ws.Range("F45").Value = "NA"
$ws.Range("G45").Value = "F, W, SP"

# Row 46
$ws.Range("C46").Value = "ME 326, ME 422."
$ws.Range("D46").Value = "NA"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "SP"

# Row 47
$ws.Range("C47").Value = "ME 329."
$ws.Range("D47").Value = "ME 318 and ME 350."
$ws.Range("E47").Value = "NA"
$ws.Range("F47").Value = "NA"
$ws.Range("G47").Value = "F, W, SP "

# Row 48
$ws.Range("C48").Value = "ME 428."
$ws.Range("D48").Value = "NA"
$ws.Range("E48").Value = "NA"
$ws.Range("F48").Value = "NA"
$ws.Range("G48").Value = "F, W, SP"

# Row 49
$ws.Range("C49").Value = "ME 429."
$ws.Range("D49").Value = "NA"
$ws.Range("E49").Value = "NA"
$ws.Range("F49").Value = "NA"
$ws.Range("G49").Value = "F, W, SP"

# Row 50
$ws.Range("C50").Value = "ME 329."
$ws.Range("D50").Value = "NA"
$ws.Range("E50").Value = "NA"
$ws.Range("F50").Value = "NA"
$ws.Range("G50").Value = "TBD"

# Row 51
$ws.Range("C51").Value = "ME 302, ME 347, and ME 350."
$ws.Range("D51").Value = "NA"
$ws.Range("E51").Value = "NA"
$ws.Range("F51").Value = "NA"
$ws.Range("G51").Value = "W"

# Row 52
$ws.Range("C52").Value = "ME 329, ME 347."
$ws.Range("D52").Value = "NA"
$ws.Range("E52").Value = "NA"
$ws.Range("F52").Value = "NA"
$ws.Range("G52").Value = "F"

# Row 53
$ws.Range("C53").Value = "ME 329, ME 347."
$ws.Range("D53").Value = "NA"
$ws.Range("E53").Value = "NA"
$ws.Range("F53").Value = "NA"
$ws.Range("G53").Value = "SP"

# Row 54
$ws.Range("C54").Value = "Junior or senior standing in the College of Engineering, ME 302, and PHYS 133."
$ws.Range("D54").Value = "NA"
$ws.Range("E54").Value = "NA"
$ws.Range("F54").Value = "NA"
$ws.Range("G54").Value = "F"

# Row 55
$ws.Range("C55").Value = "ME 302 and ME 341."
$ws.Range("D55").Value = "NA"
$ws.Range("E55").Value = "NA"
$ws.Range("F55").Value = "NA"
$ws.Range("G55").Value = "W"

# Row 56
$ws.Range("C56").Value = "ME 437 or ME 438."
$ws.Range("D56").Value = "NA"
$ws.Range("E56").Value = "NA"
$ws.Range("F56").Value = "NA"
$ws.Range("G56").Value = "SP"

# Row 57
$ws.Range("C57").Value = "ME 318, ME 329, or consent of instructor."
$ws.Range("D57").Value = "NA"
$ws.Range("E57").Value = "NA"
$ws.Range("F57").Value = "NA"
$ws.Range("G57").Value = "W"

# Row 58
$ws.Range("C58").Value = "ME 212."
$ws.Range("D58").Value = "NA"
$ws.Range("E58").Value = "NA"
$ws.Range("F58").Value = "NA"
$ws.Range("G58").Value = "F"

# Row 59
$ws.Range("C59").Value = "MATH 344, ME 303, ME 347, and ME 350."
$ws.Range("D59").Value = "NA"
$ws.Range("E59").Value = "NA"
$ws.Range("F59").Value = "NA"
$ws.Range("G59").Value = "SP"

# Row 60
$ws.Range("C60").Value = "ME 303, ME 347, and ME 350."
$ws.Range("D60").Value = "NA"
$ws.Range("E60").Value = "NA"
$ws.Range("F60").Value = "NA"
$ws.Range("G60").Value = "W"

# Row 61
$ws.Range("C61").Value = "ME 350."
$ws.Range("D61").Value = "NA"
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "ME 415."
$ws.Range("G61").Value = "W "

# Row 62
$ws.Range("C62").Value = "ME 302 or ENVE 304."
$ws.Range("D62").Value = "NA"
$ws.Range("E62").Value = "NA"
$ws.Range("F62").Value = "NA"
$ws.Range("G62").Value = "SP"

# Row 63
$ws.Range("C63").Value = "ME 350."
$ws.Range("D63").Value = "NA"
$ws.Range("E63").Value = "NA"
$ws.Range("F63").Value = "ME 359."
$ws.Range("G63").Value = "W "

# Row 64
$ws.Range("C64").Value = "ME 350."
$ws.Range("D64").Value = "NA"
$ws.Range("E64").Value = "NA"
$ws.Range("F64").Value = "ME 359."
$ws.Range("G64").Value = "SP "

# Row 65
$ws.Range("C65").Value = "ME 302, ME 347."
$ws.Range("D65").Value = "NA"
$ws.Range("E65").Value = "NA"
$ws.Range("F65").Value = "NA"
$ws.Range("G65").Value = "F"

# Row 66
$ws.Range("C66").Value = "ME 341 and ME 350."
$ws.Range("D66").Value = "NA"
$ws.Range("E66").Value = "NA"
$ws.Range("F66").Value = "NA"
$ws.Range("G66").Value = "SP"

# Row 67
$ws.Range("C67").Value = "ME 303 and ME 350."
$ws.Range("D67").Value = "NA"
$ws.Range("E67").Value = "NA"
$ws.Range("F67").Value = "NA"
$ws.Range("G67").Value = "F"

# Row 68
$ws.Range("C68").Value = "ME 456, ME 458."
$ws.Range("D68").Value = "NA"
$ws.Range("E68").Value = "NA"
$ws.Range("F68").Value = "NA"
$ws.Range("G68").Value = "W"

# Row 69
$ws.Range("C69").Value = "ME 459."
$ws.Range("D69").Value = "NA"
$ws.Range("E69").Value = "NA"
$ws.Range("F69").Value = "NA"
$ws.Range("G69").Value = "SP"

# Row 70
$ws.Range("C70").Value = "Consent of instructor."
$ws.Range("D70").Value = "NA"
$ws.Range("E70").Value = "NA"
$ws.Range("F70").Value = "NA"
$ws.Range("G70").Value = "TBD"

# Row 71
$ws.Range("C71").Value = "Consent of instructor."
$ws.Range("D71").Value = "NA"
$ws.Range("E71").Value = "NA"
$ws.Range("F71").Value = "NA"
$ws.Range("G71").Value = "TBD"

# Row 72
$ws.Range("C72").Value = "ME 329, ME 347, ME 302."
$ws.Range("D72").Value = "NA"
$ws.Range("E72").Value = "NA"
$ws.Range("F72").Value = "NA"
$ws.Range("G72").Value = "SP"

# Row 73
$ws.Range("C73").Value = "Sophomore standing and consent of instructor."
$ws.Range("D73").Value = "NA"
$ws.Range("E73").Value = "NA"
$ws.Range("F73").Value = "NA"
$ws.Range("G73").Value = "F, W, SP"

# Row 74
$ws.Range("C74").Value = "Sophomore standing and consent of instructor."
$ws.Range("D74").Value = "NA"
$ws.Range("E74").Value = "NA"
$ws.Range("F74").Value = "NA"
$ws.Range("G74").Value = "F, W, SP"

# Row 75
$ws.Range("C75").Value = "Sophomore standing and consent of instructor."
$ws.Range("D75").Value = "NA"
$ws.Range("E75").Value = "NA"
$ws.Range("F75").Value = "NA"
$ws.Range("G75").Value = "F, W, SP"

# Row 76
$ws.Range("C76").Value = "Consent of department head, graduate advisor and supervising faculty member."
$ws.Range("D76").Value = "NA"
$ws.Range("E76").Value = "NA"
$ws.Range("F76").Value = "NA"
$ws.Range("G76").Value = "F, W, SP"

# Row 77
$ws.Range("C77").Value = "Graduate standing."
$ws.Range("D77").Value = "NA"
$ws.Range("E77").Value = "NA"
$ws.Range("F77").Value = "NA"
$ws.Range("G77").Value = "TBD"

# Row 78
$ws.Range("C78").Value = "ME 501 or CE 511."
$ws.Range("D78").Value = "NA"
$ws.Range("E78").Value = "NA"
$ws.Range("F78").Value = "NA"
$ws.Range("G78").Value = "TBD"

# Row 79
$ws.Range("C79").Value = "CE/ME 404 and CE 511/ME 501 or consent of instructor."
$ws.Range("D79").Value = "NA"
$ws.Range("E79").Value = "NA"
$ws.Range("F79").Value = "NA"
$ws.Range("G79").Value = "SP"

# Row 80
$ws.Range("C80").Value = "Graduate standing or consent of instructor."
$ws.Range("D80").Value = "NA"
$ws.Range("E80").Value = "NA"
$ws.Range("F80").Value = "NA"
$ws.Range("G80").Value = "SP"

# Row 81
$ws.Range("C81").Value = "Graduate standing or consent of instructor."
$ws.Range("D81").Value = "NA"
$ws.Range("E81").Value = "NA"
$ws.Range("F81").Value = "NA"
$ws.Range("G81").Value = "F"

# Row 82
$ws.Range("C82").Value = "ME 318, graduate standing or consent of instructor."
$ws.Range("D82").Value = "NA"
$ws.Range("E82").Value = "NA"
$ws.Range("F82").Value = "NA"
$ws.Range("G82").Value = "SP"

# Row 83
$ws.Range("C83").Value = "ME 318, graduate standing or consent of instructor."
$ws.Range("D83").Value = "NA"
$ws.Range("E83").Value = "NA"
$ws.Range("F83").Value = "NA"
$ws.Range("G83").Value = "W"

# Row 84
$ws.Range("C84").Value = "ME 347, MATH 344 and graduate standing or consent of instructor."
$ws.Range("D84").Value = "NA"
$ws.Range("E84").Value = "NA"
$ws.Range("F84").Value = "NA"
$ws.Range("G84").Value = "W"

# Row 85
$ws.Range("C85").Value = "ME 303, ME 347, ME 350, and graduate standing."
$ws.Range("D85").Value = "NA"
$ws.Range("E85").Value = "NA"
$ws.Range("F85").Value = "NA"
$ws.Range("G85").Value = "TBD"

# Row 86
$ws.Range("C86").Value = "MATH 244, ME 303, ME 347, ME 350, and graduate standing."
$ws.Range("D86").Value = "NA"
$ws.Range("E86").Value = "NA"
$ws.Range("F86").Value = "NA"
$ws.Range("G86").Value = "F"

# Row 87
$ws.Range("C87").Value = "MATH 344, ME 347, ME 350, and graduate standing."
$ws.Range("D87").Value = "NA"
$ws.Range("E87").Value = "NA"
$ws.Range("F87").Value = "NA"
$ws.Range("G87").Value = "F"

# Row 88
$ws.Range("C88").Value = "MATH 344, ME 347, ME 350, and graduate standing."
$ws.Range("D88").Value = "NA"
$ws.Range("E88").Value = "NA"
$ws.Range("F88").Value = "NA"
$ws.Range("G88").Value = "W"

# Row 89
$ws.Range("C89").Value = "MATH 418, ME 347, ME 350, and graduate standing."
$ws.Range("D89").Value = "NA"
$ws.Range("E89").Value = "NA"
$ws.Range("F89").Value = "NA"
$ws.Range("G89").Value = "SP"

# Row 90
$ws.Range("C90").Value = "ME 347 or FPE 502; and ME 350."
$ws.Range("D90").Value = "NA"
$ws.Range("E90").Value = "NA"
$ws.Range("F90").Value = "NA"
$ws.Range("G90").Value = "SP"

# Row 91
$ws.Range("C91").Value = "Graduate standing in mechanical engineering program."
$ws.Range("D91").Value = "NA"
$ws.Range("E91").Value = "NA"
$ws.Range("F91").Value = "NA"
$ws.Range("G91").Value = "W"

# Row 92
$ws.Range("C92").Value = "Graduate standing or consent of instructor."
$ws.Range("D92").Value = "NA"
$ws.Range("E92").Value = "NA"
$ws.Range("F92").Value = "NA"
$ws.Range("G92").Value = "TBD"

# Row 93
$ws.Range("C93").Value = "Graduate standing of consent of instructor."
$ws.Range("D93").Value = "NA"
$ws.Range("E93").Value = "NA"
$ws.Range("F93").Value = "NA"
$ws.Range("G93").Value = "TBD"

# Row 94
$ws.Range("C94").Value = "ME 422."
$ws.Range("D94").Value = "NA"
$ws.Range("E94").Value = "NA"
$ws.Range("F94").Value = "NA"
$ws.Range("G94").Value = "TBD"

# Row 95
$ws.Range("C95").Value = "Graduate standing."
$ws.Range("D95").Value = "NA"
$ws.Range("E95").Value = "NA"
$ws.Range("F95").Value = "NA"
$ws.Range("G95").Value = "F, W, SP"
